$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "McKinnon" row (row 2) entirely - its data is not present in the
# updated table at all.
$ws.Rows.Item(2).Delete()

# Row 2 is now the former "Melbourne" row - refresh its values (time/notes
# wording changed, "Exist" flag changed from old -> new).
$ws.Cells.Item(2, 1).Value = "Melbourne"
$ws.Cells.Item(2, 2).Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Cells.Item(2, 3).Value = "28/12/2020 10:30pm-12.00am"
$ws.Cells.Item(2, 4).Value = "Case attended venue"
$ws.Cells.Item(2, 5).Value = "new"

# Row 3 is now the former "Southbank" row - replace entirely with the new
# "Moorabbin" location details.
$ws.Cells.Item(3, 1).Value = "Moorabbin"
$ws.Cells.Item(3, 2).Value = "Grape and Grain Liquor Cellars, 14/16 Station St"
$ws.Cells.Item(3, 3).Value = "21/12/20 2pm - 10pm  22/12/20 10am - 6pm  24/12/20 1pm - 10pm  28/12/20 8.05pm - 8.47pm  29/12/20 12pm - 4pm"
$ws.Cells.Item(3, 4).Value = "Case's workplace"
$ws.Cells.Item(3, 5).Value = "new"

# Resize columns to fit the new (wider/narrower) content.
$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(3).EntireColumn.AutoFit()
$ws.Columns.Item(4).EntireColumn.AutoFit()
$ws.Columns.Item(5).EntireColumn.AutoFit()

# Update the selection to full columns A:E, matching the saved view state.
$ws.Range("A1:E1048576").Select()
